# The source workbook tracks daily "Ajo" (garlic) price records for the
# "Feria Lagunitas de Puerto Montt" market. This edit adds one new weekly
# record, inserted as a new row 130 (pushing the existing rows 130-208
# down to 131-209), matching the reference diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 130; everything currently at/after
# row 130 shifts down by one (dimension grows from R208 to R209).
$ws.Rows.Item(130).Insert()

# Populate the newly inserted row 130 with the new record.
$ws.Cells.Item(130, 1).Value = 4
$ws.Cells.Item(130, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(130, 3).Value = "Los Lagos"
$ws.Cells.Item(130, 4).Value = 44582
$ws.Cells.Item(130, 5).Value = 10
$ws.Cells.Item(130, 6).Value = 100112003
$ws.Cells.Item(130, 7).Value = "Ajo"
$ws.Cells.Item(130, 8).Value = "Chilote"
$ws.Cells.Item(130, 9).Value = "Primera"
$ws.Cells.Item(130, 10).Value = 240
$ws.Cells.Item(130, 11).Value = 21000
$ws.Cells.Item(130, 12).Value = 22000
$ws.Cells.Item(130, 13).Value = 21500
$ws.Cells.Item(130, 14).Value = "`$/caja 10 kilos"
$ws.Cells.Item(130, 15).Value = "China"
$ws.Cells.Item(130, 16).Value = 2150
$ws.Cells.Item(130, 17).Value = 10
$ws.Cells.Item(130, 18).Value = "Hortaliza"
